# Modify sort and select theme sequence diagrams
#
# The only deliberate content edit in this commit is a resize/reposition
# of the big rounded-rectangle "swimlane" background shape on slide 1
# ("Rectangle 65") so that it spans the full slide width instead of being
# inset from the left edge.
#
# Target OOXML (EMU):
#   <a:off x="0" y="1752600"/>
#   <a:ext cx="9144000" cy="4419600"/>
#
# PowerPoint's Shape.Left/Top/Width/Height COM properties are expressed in
# points (1 pt = 12700 EMU), so convert accordingly:
#   x:  0       EMU ->   0  pt
#   y:  1752600 EMU -> 138  pt
#   cx: 9144000 EMU -> 720  pt
#   cy: 4419600 EMU -> 348  pt

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item(1)

$shp.Left   = 0
$shp.Top    = 138
$shp.Width  = 720
$shp.Height = 348
